$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A1: header text changed
$ws.Range("A1").Value = "distance sensor test"

# New note next to "distance sensor delay" header (row 41), regular (non-bold) style like A53
$ws.Range("A53").Copy()
$ws.Range("C41").PasteSpecial(-4122)
$ws.Range("C41").Value = "(auto kreeg stop command bij 150m. rechter colom is wanneer die stil stond)"

# Row 53 column headers renamed
$ws.Range("C53").Value = "status opvragen"
$ws.Range("E53").Value = "command versturen"

# New row 77: description under "dynamic sensor measurements" header (row 76), regular style
$ws.Range("A53").Copy()
$ws.Range("A77").PasteSpecial(-4122)
$ws.Range("A77").Value = "auto reed vanuit 5m naar 0m  toe"

# New empty (but styled, regular) cell at C85
$ws.Range("A53").Copy()
$ws.Range("C85").PasteSpecial(-4122)

# New row 91: trailing note, regular style
$ws.Range("A53").Copy()
$ws.Range("A91").PasteSpecial(-4122)
$ws.Range("A91").Value = "hangt af van obstakels in fov van sensor"

$excel.CutCopyMode = $false
